$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.804.20"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "2.249.29"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.83%  "
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.487"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").Value = "2.595.00"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "2.264.46"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.786"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "41.672.18"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.37%  "
$ws.Range("D20").Value = "0.0₃0899"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.07%  "
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0731"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").Value = "1.957.62"
$ws.Range("E43").Value = "  -2.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0280"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "71.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.08%  "
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "90.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.11%  "
